# Updates power-flow result values in Sheet1 (pl_mw) for the 380 kV case.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$dataB = @(0.3504725059204077, 0.3067962237674919, 0.2798900867041709, 0.2689041013607323, 0.2670786107846936, 0.279742011926345, 0.3354318006182098, 0.4439056410939202, 0.5231218670642193, 0.5590489308034989, 0.5726372330994423, 0.5697114963139143, 0.5601671841048983, 0.554318843659388, 0.5207716867633394, 0.5001631198233554, 0.4882994276139243, 0.4842808612495162, 0.5023579977261932, 0.562971034457064, 0.6024886078279792, 0.5814064501652751, 0.5013657418063246, 0.4146424860296634)
for ($i = 0; $i -lt $dataB.Length; $i++) {
    $ws.Cells.Item(2 + $i, 2).Value = $dataB[$i]
}

$dataC = @(0.06187472207288636, 0.0557287456371256, 0.05193563547808822, 0.05038512744557977, 0.05012738094880831, 0.05191474401610208, 0.05975968994493996, 0.07498523025151371, 0.0860705658720633, 0.09109082501568366, 0.09298853506096805, 0.09257998021055869, 0.09124701862415918, 0.09043010107218663, 0.08574201717631524, 0.08286018043256149, 0.08120051173968079, 0.08063821715454367, 0.0831671763391455, 0.0916386336742363, 0.09715562299435021, 0.0942129331281194, 0.08302839239745197, 0.07088372230663254)
for ($i = 0; $i -lt $dataC.Length; $i++) {
    $ws.Cells.Item(2 + $i, 3).Value = $dataC[$i]
}

$dataE = @(0.7791858807356533, 0.7465023004269682, 0.7267885129920586, 0.7188440264367699, 0.7175302288004985, 0.7266810103270984, 0.7678430900747486, 0.8513757701753946, 0.9144764519537318, 0.9435615044658192, 0.9546300566779991, 0.9522438143046514, 0.9444710251801354, 0.9397170875184315, 0.9125833263486953, 0.8960350271721609, 0.8865526728451982, 0.8833482576496152, 0.8977929178045372, 0.946752597939394, 0.9790692613905065, 0.9617920910171307, 0.8969980773824915, 0.8284752984282875)
for ($i = 0; $i -lt $dataE.Length; $i++) {
    $ws.Cells.Item(2 + $i, 5).Value = $dataE[$i]
}

$dataF = @(2.084801462075916, 2.059924279027712, 2.045942136159567, 2.040568928464197, 2.039696306836731, 2.045868357314873, 2.075955290454033, 2.145235300228478, 2.202443894595731, 2.229849005995732, 2.240425743223255, 2.238138994601371, 2.23071516917642, 2.22619379391017, 2.200680729225098, 2.185383246821175, 2.176714462415873, 2.173801660684092, 2.186998239914985, 2.232890321555089, 2.264043694985958, 2.24731023064362, 2.186267709266801, 2.125388402271994)
for ($i = 0; $i -lt $dataF.Length; $i++) {
    $ws.Cells.Item(2 + $i, 6).Value = $dataF[$i]
}

$dataG = @(0.3355176703624885, 0.3413374901511581, 0.3452129829686967, 0.3468680766770973, 0.3471474769070113, 0.345234997455691, 0.3374615437623412, 0.3246222282724744, 0.3166649825801784, 0.3133676216437919, 0.3121655423278682, 0.3124223577673675, 0.3132677916159494, 0.313791713475851, 0.3168869776419783, 0.3188685427063618, 0.3200386235599169, 0.3204399973077514, 0.3186544605405999, 0.3130182020352024, 0.309606056825757, 0.3114022787407649, 0.3187511509977483, 0.3278371589872791)
for ($i = 0; $i -lt $dataG.Length; $i++) {
    $ws.Cells.Item(2 + $i, 7).Value = $dataG[$i]
}

$dataH = @(0.5226100401035509, 0.5289019443049661, 0.5330191474058488, 0.5347607918027499, 0.5350538467618442, 0.5330423773111441, 0.5247267781224139, 0.5104350354312857, 0.5011633684872621, 0.4972122023250094, 0.4957543403758962, 0.4960666105640144, 0.4970914942155034, 0.4977242613360531, 0.5014269527496751, 0.5037667264434873, 0.505137590754174, 0.5056060499691739, 0.5035150567870872, 0.4967894197710194, 0.4926174254036795, 0.4948236288433705, 0.5036287565942033, 0.5140855987005608)
for ($i = 0; $i -lt $dataH.Length; $i++) {
    $ws.Cells.Item(2 + $i, 8).Value = $dataH[$i]
}

$dataI = @(0.3915410455148978, 0.3985927972270229, 0.4031865660489879, 0.4051248675279853, 0.4054507234423745, 0.4032124382973352, 0.3939176607297634, 0.3777880681614896, 0.3672200502992418, 0.3626916345092983, 0.3610170429211745, 0.3613759057855379, 0.3625530582020726, 0.3632793389189306, 0.367521615999717, 0.3701956506620014, 0.3717599423465963, 0.3722940914518347, 0.3699082768318718, 0.3622062076362269, 0.3574069451863178, 0.3599469143581562, 0.3700381145888958, 0.38192661902794)
for ($i = 0; $i -lt $dataI.Length; $i++) {
    $ws.Cells.Item(2 + $i, 9).Value = $dataI[$i]
}

$dataJ = @(0.03220732150366157, 0.03271019660061292, 0.03304176368915712, 0.03318260490973124, 0.03320633701843612, 0.03304363995398329, 0.03237597583319562, 0.03124795385230783, 0.03053018056527534, 0.03022784936530343, 0.03011685169712308, 0.03014060170789001, 0.03021864754398251, 0.03026690750081507, 0.03055042622823301, 0.0307305582177726, 0.03083644138678032, 0.03087268219696782, 0.0307111472257624, 0.03019562882661653, 0.02987904682675513, 0.03004614810897976, 0.03071991570066324, 0.03153365842561584)
for ($i = 0; $i -lt $dataJ.Length; $i++) {
    $ws.Cells.Item(2 + $i, 10).Value = $dataJ[$i]
}

$dataK = @(0.3683213754531209, 0.3214746140497766, 0.2925936336646373, 0.2807958410727167, 0.2788351262796027, 0.2924346390354913, 0.3521933543580644, 0.4684216205588996, 0.5531967362810803, 0.5916224150017513, 0.6061525108909223, 0.6030241373660203, 0.5928182382030514, 0.586564085305298, 0.5506826552964412, 0.5286343349047513, 0.5159396990398193, 0.511639303703447, 0.5309827707937416, 0.5958165318211002, 0.6380670568677544, 0.6155286217301068, 0.5299211013050638, 0.4370848806532592)
for ($i = 0; $i -lt $dataK.Length; $i++) {
    $ws.Cells.Item(2 + $i, 11).Value = $dataK[$i]
}

$dataO = @(1.650432243815047, 1.675728951735252, 1.692426282854726, 1.699523170415851, 1.700719263977774, 1.692520809773256, 1.658912521137623, 1.602267093407747, 1.56631656932089, 1.551196860365067, 1.545649329663831, 1.546836166621901, 1.550736892075591, 1.55314939071755, 1.567329557035563, 1.576345118237626, 1.581646786490978, 1.583461775681727, 1.575373371212507, 1.549586320485403, 1.533770586157246, 1.542116644596462, 1.575812328949553, 1.616597574751111)
for ($i = 0; $i -lt $dataO.Length; $i++) {
    $ws.Cells.Item(2 + $i, 15).Value = $dataO[$i]
}
